# Apply the changes described by the commit:
# "Updated arguments document, exploratory analysis, and problem model"
#
# Concrete, reproducible edits made to the "Problem Model" workbook:
#   1. Add two new workbook-level defined names (Population_Size / PopulationSize)
#      that refer to #REF! (as seen in the target OOXML).
#   2. On the "Costs" sheet, bump the window zoom level from 150% to 170%.
#   3. On the "Impact" sheet:
#        - change the selected cell from G10 to A24
#        - change the Lift assumption (B1) from 1.2 to 2.18, which cascades
#          through all of the dependent formulas on the sheet.

$wb = $excel.ActiveWorkbook

# --- 1. New defined names -------------------------------------------------
$wb.Names.Add("Population_Size", "=#REF!")
$wb.Names.Add("PopulationSize", "=#REF!")

# --- 2. Costs sheet: zoom 150% -> 170% ------------------------------------
$wsCosts = $wb.Worksheets.Item("Costs")
$wsCosts.Activate()
$excel.ActiveWindow.Zoom = 170

# --- 3. Impact sheet: selection + Lift value ------------------------------
$wsImpact = $wb.Worksheets.Item("Impact")
$wsImpact.Activate()
$wsImpact.Range("B1").Value = 2.18
$wsImpact.Range("A24").Select()

# Restore the originally active sheet/tab (Dynamics) so the workbook-level
# active tab and the per-sheet `tabSelected` flag are left exactly as they
# were before this edit (the diff does not touch sheet1 / Dynamics).
$wsDynamics = $wb.Worksheets.Item("Dynamics")
$wsDynamics.Activate()
